$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.2035556666666667
$ws.Range("H2").Value = 0.6106670000000001
$ws.Range("I2").Value = 0.006148914270823412
$ws.Range("J2").Value = 0.006148914270823412
$ws.Range("M2").Value = 3.241087666666667
$ws.Range("N2").Value = 9.723262999999999
$ws.Range("O2").Value = 0.02486257877280725
$ws.Range("P2").Value = 0.02486257877280725
$ws.Range("Q2").Value = 0.6597417607134445
$ws.Range("R2").Value = 5.937675846421
$ws.Range("S2").Value = 0.0001528778654255857
$ws.Range("T2").Value = 0.0001528778654255857
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.2035556666666667
$ws.Range("H3").Value = 0.6106670000000001
$ws.Range("I3").Value = 0.006148914270823412
$ws.Range("J3").Value = 0.006148914270823412
$ws.Range("O3").Value = 0.02096124117795788
$ws.Range("P3").Value = 0.02096124117795788
$ws.Range("Q3").Value = 0.5562176911676667
$ws.Range("R3").Value = 5.005959220509
$ws.Range("S3").Value = 0.0001288888750133166
$ws.Range("T3").Value = 0.0001288888750133166
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.2035556666666667
$ws.Range("H4").Value = 0.6106670000000001
$ws.Range("I4").Value = 0.006148914270823412
$ws.Range("J4").Value = 0.006148914270823412
$ws.Range("M4").Value = 124.3864796666667
$ws.Range("N4").Value = 373.159439
$ws.Range("O4").Value = 0.9541761800492348
$ws.Range("P4").Value = 0.9541761800492349
$ws.Range("Q4").Value = 25.31957279286812
$ws.Range("R4").Value = 227.876155135813
$ws.Range("S4").Value = 0.005867147530384509
$ws.Range("T4").Value = 0.005867147530384509
$ws.Range("I5").Value = 0.735846381812327
$ws.Range("J5").Value = 0.735846381812327
$ws.Range("M5").Value = 3.241087666666667
$ws.Range("N5").Value = 9.723262999999999
$ws.Range("O5").Value = 0.02486257877280725
$ws.Range("P5").Value = 0.02486257877280725
$ws.Range("Q5").Value = 78.95191999261232
$ws.Range("R5").Value = 710.5672799335108
$ws.Range("S5").Value = 0.01829503863249418
$ws.Range("T5").Value = 0.01829503863249418
$ws.Range("I6").Value = 0.735846381812327
$ws.Range("J6").Value = 0.735846381812327
$ws.Range("O6").Value = 0.02096124117795788
$ws.Range("P6").Value = 0.02096124117795788
$ws.Range("S6").Value = 0.01542425347909587
$ws.Range("T6").Value = 0.01542425347909587
$ws.Range("I7").Value = 0.735846381812327
$ws.Range("J7").Value = 0.735846381812327
$ws.Range("M7").Value = 124.3864796666667
$ws.Range("N7").Value = 373.159439
$ws.Range("O7").Value = 0.9541761800492348
$ws.Range("P7").Value = 0.9541761800492349
$ws.Range("Q7").Value = 3030.01720434962
$ws.Range("R7").Value = 27270.15483914658
$ws.Range("S7").Value = 0.7021270897007369
$ws.Range("T7").Value = 0.7021270897007369
$ws.Range("G8").Value = 8.541072
$ws.Range("H8").Value = 25.623216
$ws.Range("I8").Value = 0.2580047039168495
$ws.Range("J8").Value = 0.2580047039168495
$ws.Range("M8").Value = 3.241087666666667
$ws.Range("N8").Value = 9.723262999999999
$ws.Range("O8").Value = 0.02486257877280725
$ws.Range("P8").Value = 0.02486257877280725
$ws.Range("Q8").Value = 27.682363119312
$ws.Range("R8").Value = 249.141268073808
$ws.Range("S8").Value = 0.006414662274887482
$ws.Range("T8").Value = 0.006414662274887483
$ws.Range("G9").Value = 8.541072
$ws.Range("H9").Value = 25.623216
$ws.Range("I9").Value = 0.2580047039168495
$ws.Range("J9").Value = 0.2580047039168495
$ws.Range("O9").Value = 0.02096124117795788
$ws.Range("P9").Value = 0.02096124117795788
$ws.Range("Q9").Value = 23.338556109648
$ws.Range("R9").Value = 210.047004986832
$ws.Range("S9").Value = 0.005408098823848699
$ws.Range("T9").Value = 0.005408098823848699
$ws.Range("G10").Value = 8.541072
$ws.Range("H10").Value = 25.623216
$ws.Range("I10").Value = 0.2580047039168495
$ws.Range("J10").Value = 0.2580047039168495
$ws.Range("M10").Value = 124.3864796666667
$ws.Range("N10").Value = 373.159439
$ws.Range("O10").Value = 0.9541761800492348
$ws.Range("P10").Value = 0.9541761800492349
$ws.Range("Q10").Value = 1062.393878659536
$ws.Range("R10").Value = 9561.544907935824
$ws.Range("S10").Value = 0.2461819428181133
$ws.Range("T10").Value = 0.2461819428181134
